# Update latest output (run 184)

$wb = $excel.ActiveWorkbook

# --- Sheet "Schedule" (E2, F2) ---
$wsSchedule = $wb.Worksheets.Item("Schedule")
$wsSchedule.Range("E2").Value = 1309.133241
$wsSchedule.Range("F2").Value = 21.64572157738095

# --- Sheet "Detailed" ---
$wsDetailed = $wb.Worksheets.Item("Detailed")

$wsDetailed.Range("B8").Value = 56.97994
$wsDetailed.Range("B9").Value = 56.98
$wsDetailed.Range("C10").Value = "historical"
$wsDetailed.Range("B13").Value = 64.8901
$wsDetailed.Range("B17").Value = 35.88
$wsDetailed.Range("B18").Value = 25.17119
$wsDetailed.Range("B19").Value = 34.16904
$wsDetailed.Range("B21").Value = 1.46051
$wsDetailed.Range("B22").Value = 36.06
$wsDetailed.Range("B25").Value = 34.01
$wsDetailed.Range("B26").Value = 36.06
$wsDetailed.Range("B33").Value = 0.7
$wsDetailed.Range("B34").Value = 0.06
$wsDetailed.Range("B35").Value = -1.589
$wsDetailed.Range("B36").Value = -2.15325
$wsDetailed.Range("B37").Value = 4.80939
$wsDetailed.Range("B38").Value = 44.27664
$wsDetailed.Range("B39").Value = 45.44136
$wsDetailed.Range("B40").Value = 57.03871
$wsDetailed.Range("B41").Value = 58.88652
$wsDetailed.Range("B42").Value = 59.38951
